$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 37
$ws.Range("B3").Value = 44
$ws.Range("B4").Value = 30
$ws.Range("B5").Value = 17
$ws.Range("B6").Value = 52
